$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '20.571.71'
$ws.Cells.Item(2, 5).Value = '  +0.24%  '

$ws.Cells.Item(3, 4).Value = '1.478.58'
$ws.Cells.Item(3, 5).Value = '  +0.54%  '

$ws.Cells.Item(4, 5).Value = '  +0.19%  '

$ws.Cells.Item(5, 4).Value = '''0.9733'
$ws.Cells.Item(5, 5).Value = '  +1.45%  '

$ws.Cells.Item(6, 4).Value = '''279.31'
$ws.Cells.Item(6, 5).Value = '  -0.92%  '

$ws.Cells.Item(7, 4).Value = '''0.3660'
$ws.Cells.Item(7, 5).Value = '  -1.06%  '

$ws.Cells.Item(8, 4).Value = '''0.3077'
$ws.Cells.Item(8, 5).Value = '  -3.35%  '

$ws.Cells.Item(9, 4).Value = '''40.04'
$ws.Cells.Item(9, 5).Value = '  -4.25%  '

$ws.Cells.Item(10, 4).Value = '''1.061'
$ws.Cells.Item(10, 5).Value = '  +0.56%  '

$ws.Cells.Item(11, 4).Value = '''0.06669'
$ws.Cells.Item(11, 5).Value = '  -0.23%  '

$ws.Cells.Item(12, 5).Value = '  +0.01%  '

$ws.Cells.Item(13, 4).Value = '''5.514'
$ws.Cells.Item(13, 5).Value = '  -1.87%  '

$ws.Cells.Item(14, 4).Value = '''18.05'
$ws.Cells.Item(14, 5).Value = '  -1.14%  '

$ws.Cells.Item(15, 4).Value = '''6.211'
$ws.Cells.Item(15, 5).Value = '  -0.91%  '

$ws.Cells.Item(16, 4).Value = '''0.9746'
$ws.Cells.Item(16, 5).Value = '  +1.86%  '

$ws.Cells.Item(17, 4).Value = '''0.00001028'
$ws.Cells.Item(17, 5).Value = '  -0.70%  '

$ws.Cells.Item(18, 4).Value = '1.477.61'
$ws.Cells.Item(18, 5).Value = '  +0.25%  '

$ws.Cells.Item(19, 4).Value = '''0.05923'
$ws.Cells.Item(19, 5).Value = '  +4.42%  '

$ws.Cells.Item(20, 4).Value = '''69.43'
$ws.Cells.Item(20, 5).Value = '  -3.95%  '

$ws.Cells.Item(21, 4).Value = '''5.484'
$ws.Cells.Item(21, 5).Value = '  -3.50%  '

$ws.Cells.Item(22, 4).Value = '''14.50'
$ws.Cells.Item(22, 5).Value = '  -1.45%  '

$ws.Cells.Item(23, 4).Value = '''11.05'
$ws.Cells.Item(23, 5).Value = '  -1.51%  '

$ws.Cells.Item(24, 4).Value = '''2.259'
$ws.Cells.Item(24, 5).Value = '  -0.56%  '

$ws.Cells.Item(25, 4).Value = '20.622.75'
$ws.Cells.Item(25, 5).Value = '  -0.30%  '

$ws.Cells.Item(26, 4).Value = '''142.20'
$ws.Cells.Item(26, 5).Value = '  +3.17%  '

$ws.Cells.Item(27, 4).Value = '''2.128'
$ws.Cells.Item(27, 5).Value = '  -7.31%  '

$ws.Cells.Item(28, 4).Value = '''17.26'
$ws.Cells.Item(28, 5).Value = '  -1.74%  '

$ws.Cells.Item(29, 4).Value = '1.637.62'
$ws.Cells.Item(29, 5).Value = '  +0.02%  '

$ws.Cells.Item(30, 4).Value = '''114.12'
$ws.Cells.Item(30, 5).Value = '  +0.14%  '

$ws.Cells.Item(31, 4).Value = '''3.946'
$ws.Cells.Item(31, 5).Value = '  -0.14%  '

$ws.Cells.Item(32, 4).Value = '''5.012'
$ws.Cells.Item(32, 5).Value = '  -5.87%  '

$ws.Cells.Item(33, 4).Value = '''0.8185'
$ws.Cells.Item(33, 5).Value = '  -2.06%  '

$ws.Cells.Item(34, 4).Value = '''0.07997'
$ws.Cells.Item(34, 5).Value = '  +1.91%  '

$ws.Cells.Item(35, 4).Value = '''1.540'
$ws.Cells.Item(35, 5).Value = '  -5.01%  '

$ws.Cells.Item(36, 4).Value = '''1.202'
$ws.Cells.Item(36, 5).Value = '  +7.17%  '

$ws.Cells.Item(37, 4).Value = '''0.05779'
$ws.Cells.Item(37, 5).Value = '  -4.36%  '

$ws.Cells.Item(38, 4).Value = '''4.715'
$ws.Cells.Item(38, 5).Value = '  -4.10%  '

$ws.Cells.Item(39, 4).Value = '''0.9741'
$ws.Cells.Item(39, 5).Value = '  +0.62%  '

$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(40, 4).Value = '''0.02043'
$ws.Cells.Item(40, 5).Value = '  -1.46%  '

$ws.Cells.Item(41, 2).Value = 'FraxShare'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(41, 4).Value = '''7.651'
$ws.Cells.Item(41, 5).Value = '  +3.87%  '

$ws.Cells.Item(42, 5).Value = '  -1.79%  '

$ws.Cells.Item(43, 4).Value = '''0.1886'
$ws.Cells.Item(43, 5).Value = '  -0.16%  '

$ws.Cells.Item(44, 4).Value = '''0.5298'
$ws.Cells.Item(44, 5).Value = '  -2.37%  '

$ws.Cells.Item(45, 4).Value = '''3.532'
$ws.Cells.Item(45, 5).Value = '  -1.75%  '

$ws.Cells.Item(46, 4).Value = '''12.13'
$ws.Cells.Item(46, 5).Value = '  -3.12%  '

$ws.Cells.Item(47, 4).Value = '''118.37'
$ws.Cells.Item(47, 5).Value = '  -3.23%  '

$ws.Cells.Item(48, 4).Value = '''0.5199'
$ws.Cells.Item(48, 5).Value = '  -2.61%  '

$ws.Cells.Item(49, 4).Value = '''1.805'
$ws.Cells.Item(49, 5).Value = '  -1.58%  '

$ws.Cells.Item(50, 4).Value = '''0.06486'
$ws.Cells.Item(50, 5).Value = '  +0.71%  '

$ws.Cells.Item(51, 5).Value = '  +0.17%  '
